# Risk Assessment workbook update
# Adds two new risk rows ("Feature Creep" and inconsistent work hours / time
# constraints risks) to the risk register table, extends the parallel
# risk-score legend table to match, and tidies up column widths / selection
# the way Excel would after this kind of data-entry session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-blank legend cells for the existing row 6 ---
$ws.Range("I6").Value = "High"
$ws.Range("J6").Value = "Low"
$ws.Range("K6").Value = "Low"
$ws.Range("L6").Value = 7

# --- New risk #6: Feature Creep ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Increased scope leading to loss of focus on core user requirements"
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Focus on core functionality defined in user requirements before continuing to add functionality"
$ws.Range("G7").Value = "Feature Creep"
$ws.Range("I7").Value = "High"
$ws.Range("J7").Value = "Low"
$ws.Range("K7").Value = "Low"
$ws.Range("L7").Value = 6

# --- New risk #7: Lack of consistency in work hours ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Time constraints leading to user requirements not being fully met"
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Ensure work is properly structured, using tools such as the Gantt chart"
$ws.Range("G8").Value = "Lack of consistency in work hours"
$ws.Range("I8").Value = "High"
$ws.Range("J8").Value = "Low"
$ws.Range("K8").Value = "Low"
$ws.Range("L8").Value = 5

# --- Grow the risk-score legend table (Table3) so it covers the 3 new rows ---
$legend = $ws.ListObjects.Item("Table3")
$legend.Resize($ws.Range("I1:L8"))

# --- Column widths widened to fit the newly-entered, longer text ---
$ws.Columns.Item(2).ColumnWidth = 61.666666666666664
$ws.Columns.Item(6).ColumnWidth = 126.33333333333333
$ws.Columns.Item(7).ColumnWidth = 41.666666666666664

# --- Leave the selection where it ended up after the data entry ---
$ws.Range("M11").Select() | Out-Null
